# Update countries & provincias Spain
# Refreshes the COVID-19 stats table on sheet "Pais": updated case counts
# for several countries change the descending sort order by "Casos
# totales" (column B), so some adjacent rows swap country names while
# their own numbers roll to whichever row they now occupy. Also bumps
# the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

function Set-Row($Row, $Country, $B, $C, $D, $E, $F, $G, $H) {
    $ws.Range("A$Row").Value = $Country
    $ws.Range("B$Row").Value = $B
    $ws.Range("C$Row").Value = $C
    $ws.Range("D$Row").Value = $D
    $ws.Range("E$Row").Value = $E
    $ws.Range("F$Row").Value = $F
    $ws.Range("G$Row").Value = $G
    $ws.Range("H$Row").Value = $H
}

# Estados Unidos (row 4)
Set-Row 4 "Estados Unidos" 2226527 18127 908542 1198199 0 654 119786

# India (row 7)
Set-Row 7 "India" 367264 13103 194438 160564 0 341 12262

# Peru overtakes Italia (rows 10/11 swap)
Set-Row 10 "Peru" 240908 3752 128622 105029 0 201 7257
Set-Row 11 "Italia" 237828 328 179455 23925 0 43 34448

# Canada (row 20)
Set-Row 20 "Canada" 99786 319 62000 29532 0 41 8254

# Egipto overtakes Paises Bajos & Ecuador (rows 29/30/31 rotate)
Set-Row 29 "Egipto" 49219 1363 12730 34639 0 84 1850
Set-Row 30 "Paises Bajos" 49204 117 0 0 0 4 6074
Set-Row 31 "Ecuador" 48490 547 23881 20602 0 37 4007

# Mauritania overtakes Maldivas (rows 101/102 swap)
Set-Row 101 "Mauritania" 2223 166 427 1701 0 2 95
Set-Row 102 "Maldivas" 2120 26 1677 435 0 0 8

# Costa Rica overtakes Nicaragua & Islandia (rows 106/107/108 rotate)
Set-Row 106 "Costa Rica" 1871 75 899 960 0 0 12
Set-Row 107 "Nicaragua" 1823 0 1238 521 0 0 64
Set-Row 108 "Islandia" 1815 3 1797 8 0 0 10

# Zambia (row 119)
Set-Row 119 "Zambia" 1412 7 1142 259 0 0 11

# Cabo Verde (row 136)
Set-Row 136 "Cabo Verde" 792 11 361 424 0 0 7

# Santo Tome y Principe (row 140)
Set-Row 140 "Santo Tome y Principe" 683 12 188 483 0 0 12

# Siria (row 167)
Set-Row 167 "Siria" 178 1 78 93 0 1 7

# Groenlandia overtakes Islas Malvinas (rows 206/207 swap, numbers tied)
Set-Row 206 "Groenlandia" 13 0 13 0 0 0 0
Set-Row 207 "Islas Malvinas" 13 0 13 0 0 0 0

# Islas Turcas y Caicos overtakes Santa Sede (rows 208/209 swap)
Set-Row 208 "Islas Turcas y Caicos" 12 0 11 0 0 0 1
Set-Row 209 "Santa Sede" 12 0 12 0 0 0 0

# Bump the "last updated" timestamp in row 1
$ws.Range("A1").Value = "Datos actualizados a 17 de Junio de 2020 a las 22:40"
